$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.843.04'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.736.79'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '232.18'
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '0.5179'
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("D8").Value = '0.2750'
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("D9").Value = '39.25'
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("D10").Value = '0.06124'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").Value = '1.733.28'
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D12").Value = '0.07053'
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("D13").Value = '15.18'
$ws.Range("E13").Value = '  -1.88%  '
$ws.Range("D14").Value = '0.6370'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '4.505'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '76.82'
$ws.Range("E16").Value = '  -1.55%  '
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '25.824.54'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '11.46'
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").Value = '0.000006629'
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").Value = '1.958.28'
$ws.Range("E22").Value = '  -1.46%  '
$ws.Range("D23").Value = '4.135'
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").Value = '8.723'
$ws.Range("E24").Value = '  +5.26%  '
$ws.Range("D25").Value = '5.137'
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("D26").Value = '139.35'
$ws.Range("E26").Value = '  +2.30%  '
$ws.Range("D27").Value = '1.511'
$ws.Range("E27").Value = '  +2.11%  '
$ws.Range("D28").Value = '14.99'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").Value = '1.777'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '101.88'
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").Value = '0.08296'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").Value = '3.681'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '3.477'
$ws.Range("E33").Value = '  +2.32%  '
$ws.Range("D34").Value = '0.04503'
$ws.Range("E34").Value = '  +2.55%  '
$ws.Range("D35").Value = '2.615'
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("D36").Value = '0.9732'
$ws.Range("E36").Value = '  -2.71%  '
$ws.Range("D37").Value = '0.6114'
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("D38").Value = '2.652'
$ws.Range("E38").Value = '  -2.00%  '
$ws.Range("D39").Value = '0.01579'
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").Value = '1.942'
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = '0.9996'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '100.41'
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("D43").Value = '0.3823'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").Value = '0.7241'
$ws.Range("E44").Value = '  -3.36%  '
$ws.Range("D45").Value = '4.986'
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("D46").Value = '0.05381'
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("D47").Value = '0.1123'
$ws.Range("E47").Value = '  +1.64%  '
$ws.Range("D48").Value = '6.230'
$ws.Range("E48").Value = '  +4.04%  '
$ws.Range("D49").Value = '52.92'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").Value = '29.95'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").Value = '7.622'
$ws.Range("E51").Value = '  +2.61%  '
